$d = $word.ActiveDocument

# Update Sprint No. from "1" to "2".
# Scope the Find/Replace tightly to the "Sprint No." value cell (rather than
# the whole document, since "1" alone is not unique) by rebuilding an
# explicit document Range from the cell's own Start/End bounds -- using
# Cell.Range directly here causes Find to wander outside the cell.
$cell = $d.Tables.Item(1).Cell(2, 4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("1", $true, $true, $false, $false, $false, $true, 0, $false, "2", 2)

# Update Review Date from "02/09/18" to "02/21/18". This string is unique
# across the document (the three visually-merged cells share one run), so a
# document-wide replace is safe.
$d.Content.Find.Execute("02/09/18", $true, $true, $false, $false, $false, $true, 1, $false, "02/21/18", 2)
